$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '30.448.17'
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.874.62'
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = '  -0.79%  '
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '0.9992'
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = '  -0.14%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '245.92'
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = '  -1.00%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.9994'
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = '  -0.03%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.4725'
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = '  -0.23%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.2871'
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = '  -2.00%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.06507'
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = '  -0.44%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '21.92'
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = '  -0.47%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '100.68'
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = '  +3.66%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.07808'
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = '  +0.12%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '1.872.76'
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = '  -0.85%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '0.7295'
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = '  -0.93%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '5.171'
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = '  -1.60%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '284.49'
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = '  +0.47%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '30.423.02'
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = '  -1.39%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '13.10'
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = '  -0.82%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.9997'
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = '  -0.04%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '0.000007492'
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = '  -0.76%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '2.117.31'
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = '  -1.02%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '5.329'
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = '  -0.04%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '0.9991'
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = '  -0.04%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '6.338'
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = '  +1.12%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '9.052'
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = '  -2.05%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '161.90'
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = '  -1.35%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '19.01'
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = '  +0.30%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '1.899'
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = '  -1.38%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '0.09690'
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = '  -0.31%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '1.323'
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = '  -1.36%  '
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = '  -0.36%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '4.234'
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = '  -1.63%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '4.155'
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = '  -1.16%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.04820'
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = '  -0.64%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '1.126'
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = '  -0.08%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.6919'
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = '  -1.02%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '2.744'
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = '  +0.91%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.01901'
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = '  +0.01%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '2.839'
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = '  +1.13%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '76.02'
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = '  -0.15%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '6.307'
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = '  -1.13%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '1.960'
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = '  -2.94%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.4223'
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = '  -1.07%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.9994'
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.8262'
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = '  -1.04%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '101.02'
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = '  -0.20%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '9.756'
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = '  +2.73%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '7.024'
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '35.01'
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = '  -1.85%  '
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = '  -0.05%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '883.98'
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = '  -3.86%  '
